# "forgot 2 more sheets"
#
# The workbook originally has a single sheet ("Sheet1") holding two stacked
# mini-tables (rows 1-9 and rows 10-17). The edit:
#   1. Renames that sheet to "10".
#   2. Duplicates it twice more, named "100" and "1000", each an exact copy
#      of the original data/formulas (so sharedStrings usage triples).
#   3. Leaves "1000" as the active sheet/tab (activeTab = index 2).
#   4. Re-points each sheet's selection:
#        "10"   -> whole used range selected, anchored near the end (G17)
#        "100"  -> whole used range selected
#        "1000" -> single cell F18 selected (just past the used range) and
#                  this is the tab Excel had showing (tabSelected).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename the existing sheet -----------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "10"

# --- Sheet 2: duplicate sheet 1 right after it, then rename ---------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "100"

# --- Sheet 3: duplicate sheet 2 right after it, then rename ---------------
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "1000"

# --- Selections / active tab ----------------------------------------------
# "10": select the full table, A1:G17
$ws1.Select()
$ws1.Range("A1:G17").Select()

# "100": select the full table, A1:G17
$ws2.Select()
$ws2.Range("A1:G17").Select()

# "1000": select a single cell just below the table and leave this sheet
# active, matching the saved view of the workbook.
$ws3.Select()
$ws3.Range("F18").Select()
